$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D for the new "credit notes" column
$ws.Columns.Item(4).ColumnWidth = 40

# Set values in the exact order the original author typed them, so that
# new shared-string table entries land at the same indices as the target.

# 1) D4 -> new shared string #43
$ws.Range("D4").Value = "Didn't attend the lab, didn't show TA the result. No longer on the D2L namelist."

# 2) D8 -> new shared string #44
$ws.Range("D8").Value = 'Attended the lab, but didn''t show TA the result. Responded on Piazza that he might use a "freebie" for this time.'

# 3) B4 -> new shared string #45 ("No")
$ws.Range("B4").Value = "No"
$ws.Range("C4").Value = 0

# 4) D22 -> new shared string #46
$ws.Range("D22").Value = "Attended the lab, but couldn't finish the two programs."

# 5) D24 -> new shared string #47
$ws.Range("D24").Value = "Delayed submission on Monday."

# 6) D35 -> new shared string #48 (trailing space preserved)
$ws.Range("D35").Value = "Attended the lab, but didn't show TA the result. "

# Remaining cells reuse existing / already-created shared strings.
$ws.Range("B8").Value = "Yes"
$ws.Range("C8").Value = 0

$ws.Range("B17").Value = "No"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = "Didn't attend the lab, didn't show TA the result. No longer on the D2L namelist."

$ws.Range("B20").Value = "No"
$ws.Range("C20").Value = 0

$ws.Range("B22").Value = "Yes"
$ws.Range("C22").Value = 0

$ws.Range("B23").Value = "No"
$ws.Range("C23").Value = 0

$ws.Range("B35").Value = "Yes"
$ws.Range("C35").Value = 0

# Update view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("E46").Select()
